$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-based values for rows 2-7 (existing FAPs -> * rows) and
# newly added rows 8-13 (new MuSCs -> * rows), per updated script output.

$data = @(
    @("FAPs","Il13","Il4ra","ECs",1,0.3333333333333333,0.1174156666666667,0.352247,0.6540874079906115,0.7393359457808691,2,1,46.9028855,93.80577099999999,0.2248968021537497,0.1641136908932902,5.507133569572832,33.042801417437,0.1471021663861235,0.1213351508721799),
    @("FAPs","Il13","Il4ra","FAPs",1,0.3333333333333333,0.1174156666666667,0.352247,0.6540874079906115,0.7393359457808691,3,1,31.834655,95.50396499999999,0.1526454509321794,0.1670846902488932,3.737887239928332,33.640985159355,0.09984346734178731,0.123531717490669),
    @("FAPs","Il13","Il4ra","Inflammatory-Mac",1,0.3333333333333333,0.1174156666666667,0.352247,0.6540874079906115,0.7393359457808691,3,1,40.94164133333334,122.824924,0.1963129584168644,0.214882852051052,4.807190111580889,43.264711004228,0.1284058341258555,0.1588706166532551),
    @("FAPs","Il13","Il4ra","MuSCs",1,0.3333333333333333,0.1174156666666667,0.352247,0.6540874079906115,0.7393359457808691,2,1,7.165723,14.331446,0.03435925467356532,0.02507294033005582,0.8413681431936666,5.048208859161999,0.02247395582992165,0.01853732605242912),
    @("FAPs","Il13","Il4ra","Neutrophils",1,0.3333333333333333,0.1174156666666667,0.352247,0.6540874079906115,0.7393359457808691,3,1,49.85523499999999,149.565705,0.2390531585124693,0.2616659893836066,5.853785654348332,52.68407088913499,0.1563616608233898,0.1934590717396156),
    @("FAPs","Il13","Il4ra","Resolving-Mac",1,0.3333333333333333,0.1174156666666667,0.352247,0.6540874079906115,0.7393359457808691,3,1,31.85278333333333,95.55834999999999,0.1527323753111719,0.1671798370931021,3.740015790272222,33.66014211245,0.09990032348353371,0.1236020629727203),
    @("MuSCs","Il13","Il4ra","ECs",1,0.5,0.062095,0.12419,0.3459125920093885,0.2606640542191307,2,1,46.9028855,93.80577099999999,0.2248968021537497,0.1641136908932902,2.9124346751225,11.64973870049,0.0777946357676262,0.04277854002111027),
    @("MuSCs","Il13","Il4ra","FAPs",1,0.5,0.062095,0.12419,0.3459125920093885,0.2606640542191307,3,1,31.834655,95.50396499999999,0.1526454509321794,0.1670846902488932,1.976772902225,11.86063741335,0.05280198359039211,0.04355297275822417),
    @("MuSCs","Il13","Il4ra","Inflammatory-Mac",1,0.5,0.062095,0.12419,0.3459125920093885,0.2606640542191307,3,1,40.94164133333334,122.824924,0.1963129584168644,0.214882852051052,2.542271218593334,15.25362731156,0.06790712429100887,0.05601223539779688),
    @("MuSCs","Il13","Il4ra","MuSCs",1,0.5,0.062095,0.12419,0.3459125920093885,0.2606640542191307,2,1,7.165723,14.331446,0.03435925467356532,0.02507294033005582,0.444955569685,1.77982227874,0.01188529884364368,0.006535614277626701),
    @("MuSCs","Il13","Il4ra","Neutrophils",1,0.5,0.062095,0.12419,0.3459125920093885,0.2606640542191307,3,1,49.85523499999999,149.565705,0.2390531585124693,0.2616659893836066,3.095760817325,18.57456490395,0.08269149768907949,0.06820691764399092),
    @("MuSCs","Il13","Il4ra","Resolving-Mac",1,0.5,0.062095,0.12419,0.3459125920093885,0.2606640542191307,3,1,31.85278333333333,95.55834999999999,0.1527323753111719,0.1671798370931021,1.977898581083333,11.8673914865,0.05283205182763822,0.04357777412038181)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
